$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.12720079081461
$ws.Range("C2").Value = 7.389495610375701
$ws.Range("D2").Value = 11.64418014807944
$ws.Range("F2").Value = 27.40736469037398
$ws.Range("G2").Value = 24.1331895144116
$ws.Range("H2").Value = 13.06662984721089
$ws.Range("I2").Value = 18.15780059896577
$ws.Range("J2").Value = 11.22194593619614
$ws.Range("M2").Value = 16.44688816836023
$ws.Range("N2").Value = 16.81516471214711
$ws.Range("O2").Value = 19.27898093677833

$ws.Range("B3").Value = 10.60478074335579
$ws.Range("C3").Value = 6.993795040531849
$ws.Range("D3").Value = 11.62785071394849
$ws.Range("F3").Value = 27.42860695487765
$ws.Range("G3").Value = 24.10620001714309
$ws.Range("H3").Value = 13.10498417423351
$ws.Range("I3").Value = 18.2494845832875
$ws.Range("J3").Value = 11.24969429818938
$ws.Range("M3").Value = 16.28293413478651
$ws.Range("N3").Value = 16.85101650653628
$ws.Range("O3").Value = 19.32683003109261

$ws.Range("B4").Value = 10.27091962452397
$ws.Range("C4").Value = 6.737780383668395
$ws.Range("D4").Value = 11.61978239228317
$ws.Range("F4").Value = 27.44883482000672
$ws.Range("G4").Value = 24.09862162478185
$ws.Range("H4").Value = 13.13077493048476
$ws.Range("I4").Value = 18.30957062600028
$ws.Range("J4").Value = 11.26836078033444
$ws.Range("M4").Value = 16.18378000705099
$ws.Range("N4").Value = 16.87468088012544
$ws.Range("O4").Value = 19.36082308212669

$ws.Range("B5").Value = 10.13174217887189
$ws.Range("C5").Value = 6.63022875097079
$ws.Range("D5").Value = 11.61698971450644
$ws.Range("F5").Value = 27.45888277783808
$ws.Range("G5").Value = 24.0977957766988
$ws.Range("H5").Value = 13.14184804651885
$ws.Range("I5").Value = 18.3350090767698
$ws.Range("J5").Value = 11.27637709800094
$ws.Range("M5").Value = 16.14379333857412
$ws.Range("N5").Value = 16.88474019285263
$ws.Range("O5").Value = 19.37583291761479

$ws.Range("B6").Value = 10.10844840872721
$ws.Range("C6").Value = 6.612177151178381
$ws.Range("D6").Value = 11.616555980661
$ws.Range("F6").Value = 27.46066017651241
$ws.Range("G6").Value = 24.09779526685194
$ws.Range("H6").Value = 13.14372073251805
$ws.Range("I6").Value = 18.3392906496041
$ws.Range("J6").Value = 11.27773293794441
$ws.Range("M6").Value = 16.13718003126555
$ws.Range("N6").Value = 16.88643567267202
$ws.Range("O6").Value = 19.37839510875386

$ws.Range("B7").Value = 10.26905504889816
$ws.Range("C7").Value = 6.736342870214513
$ws.Range("D7").Value = 11.61974272053536
$ws.Range("F7").Value = 27.44896302523173
$ws.Range("G7").Value = 24.09860132734798
$ws.Range("H7").Value = 13.13092198662075
$ws.Range("I7").Value = 18.30990984029072
$ws.Range("J7").Value = 11.26846723287996
$ws.Range("M7").Value = 16.18323898402602
$ws.Range("N7").Value = 16.87481485871959
$ws.Range("O7").Value = 19.36102082756246

$ws.Range("B8").Value = 10.94987970007704
$ws.Range("C8").Value = 7.255798672447268
$ws.Range("D8").Value = 11.63814518679444
$ws.Range("F8").Value = 27.41319696702788
$ws.Range("G8").Value = 24.12201835341789
$ws.Range("H8").Value = 13.07938908345836
$ws.Range("I8").Value = 18.18862592889294
$ws.Range("J8").Value = 11.2311754906472
$ws.Range("M8").Value = 16.39006592140544
$ws.Range("N8").Value = 16.82718410179222
$ws.Range("O8").Value = 19.29452021020377

$ws.Range("B9").Value = 12.17500887183204
$ws.Range("C9").Value = 8.168876846675474
$ws.Range("D9").Value = 11.68962962009411
$ws.Range("F9").Value = 27.40011090401245
$ws.Range("G9").Value = 24.23913301916411
$ws.Range("H9").Value = 12.99613023855865
$ws.Range("I9").Value = 17.98089905676372
$ws.Range("J9").Value = 11.17097157113754
$ws.Range("M9").Value = 16.80594192072392
$ws.Range("N9").Value = 16.7468524197598
$ws.Range("O9").Value = 19.20083052657372

$ws.Range("B10").Value = 13.00140848277691
$ws.Range("C10").Value = 8.773327179621058
$ws.Range("D10").Value = 11.73663751482245
$ws.Range("F10").Value = 27.42527472279616
$ws.Range("G10").Value = 24.36818336804511
$ws.Range("H10").Value = 12.94582961552377
$ws.Range("I10").Value = 17.84667302135396
$ws.Range("J10").Value = 11.13461981572178
$ws.Range("M10").Value = 17.11542180135218
$ws.Range("N10").Value = 16.69576186503371
$ws.Range("O10").Value = 19.15451883115161

$ws.Range("B11").Value = 13.36030544496446
$ws.Range("C11").Value = 9.033612250373976
$ws.Range("D11").Value = 11.7599664299445
$ws.Range("F11").Value = 27.44425255221377
$ws.Range("G11").Value = 24.43608970655631
$ws.Range("H11").Value = 12.92531068216898
$ws.Range("I11").Value = 17.78961198505099
$ws.Range("J11").Value = 11.11979342796482
$ws.Range("M11").Value = 17.25656862700702
$ws.Range("N11").Value = 16.6742329188262
$ws.Range("O11").Value = 19.13836458229039

$ws.Range("B12").Value = 13.49369224771396
$ws.Range("C12").Value = 9.13004781991258
$ws.Range("D12").Value = 11.76907537912176
$ws.Range("F12").Value = 27.45251820007477
$ws.Range("G12").Value = 24.46311116314333
$ws.Range("H12").Value = 12.91788079672308
$ws.Range("I12").Value = 17.76858042197637
$ws.Range("J12").Value = 11.11442496206857
$ws.Range("M12").Value = 17.31002796568656
$ws.Range("N12").Value = 16.66632608518435
$ws.Range("O12").Value = 19.13295530337558

$ws.Range("B13").Value = 13.46507806070223
$ws.Range("C13").Value = 9.109373600438147
$ws.Range("D13").Value = 11.76710146297644
$ws.Range("F13").Value = 27.45069012043845
$ws.Range("G13").Value = 24.45723377466799
$ws.Range("H13").Value = 12.9194658197292
$ws.Range("I13").Value = 17.77308429925145
$ws.Range("J13").Value = 11.11557021745363
$ws.Range("M13").Value = 17.29851475559353
$ws.Range("N13").Value = 16.66801804355038
$ws.Range("O13").Value = 19.13408878062587

$ws.Range("B14").Value = 13.37133011276505
$ws.Range("C14").Value = 9.041588809554192
$ws.Range("D14").Value = 11.7607103502933
$ws.Range("F14").Value = 27.44491096948766
$ws.Range("G14").Value = 24.43828669001617
$ws.Range("H14").Value = 12.92469260104406
$ws.Range("I14").Value = 17.78787015074507
$ws.Range("J14").Value = 11.11934683137387
$ws.Range("M14").Value = 17.26096679566851
$ws.Range("N14").Value = 16.67357749744937
$ws.Range("O14").Value = 19.13790536132674

$ws.Range("B15").Value = 13.31357678519218
$ws.Range("C15").Value = 8.999791003855469
$ws.Range("D15").Value = 11.75683124621954
$ws.Range("F15").Value = 27.4415114807768
$ws.Range("G15").Value = 24.42685071145308
$ws.Range("H15").Value = 12.92793846901764
$ws.Range("I15").Value = 17.79700198675794
$ws.Range("J15").Value = 11.12169214924558
$ws.Range("M15").Value = 17.23796766150607
$ws.Range("N15").Value = 16.67701480675693
$ws.Range("O15").Value = 19.14033536580727

$ws.Range("B16").Value = 12.97760464873539
$ws.Range("C16").Value = 8.756019894753539
$ws.Range("D16").Value = 11.73515164071609
$ws.Range("F16").Value = 27.42418571816153
$ws.Range("G16").Value = 24.36392934054278
$ws.Range("H16").Value = 12.94721816053137
$ws.Range("I16").Value = 17.85048266297067
$ws.Range("J16").Value = 11.13562317260376
$ws.Range("M16").Value = 17.10620121735852
$ws.Range("N16").Value = 16.69720324464404
$ws.Range("O16").Value = 19.15567354869232

$ws.Range("B17").Value = 12.76707946225773
$ws.Range("C17").Value = 8.602699794829375
$ws.Range("D17").Value = 11.72234664885523
$ws.Range("F17").Value = 27.41548339225243
$ws.Range("G17").Value = 24.32767489718803
$ws.Range("H17").Value = 12.95965114825376
$ws.Range("I17").Value = 17.88431638488784
$ws.Range("J17").Value = 11.14460746367238
$ws.Range("M17").Value = 17.02542961621773
$ws.Range("N17").Value = 16.71002640414524
$ws.Range("O17").Value = 19.1663425641225

$ws.Range("B18").Value = 12.6443907157478
$ws.Range("C18").Value = 8.51313286423367
$ws.Range("D18").Value = 11.71516491949597
$ws.Range("F18").Value = 27.4111871613445
$ws.Range("G18").Value = 24.30768945192238
$ws.Range("H18").Value = 12.96702467654844
$ws.Range("I18").Value = 17.90415298515103
$ws.Range("J18").Value = 11.1499359619656
$ws.Range("M18").Value = 16.97900857492306
$ws.Range("N18").Value = 16.71756315686195
$ws.Range("O18").Value = 19.17294148320279

$ws.Range("B19").Value = 12.60257790877224
$ws.Range("C19").Value = 8.482570414787762
$ws.Range("D19").Value = 11.71276494532854
$ws.Range("F19").Value = 27.40985441920543
$ws.Range("G19").Value = 24.3010721103664
$ws.Range("H19").Value = 12.9695594122428
$ws.Range("I19").Value = 17.91093391935095
$ws.Range("J19").Value = 11.15176774546259
$ws.Range("M19").Value = 16.96329879261177
$ws.Range("N19").Value = 16.72014267624058
$ws.Range("O19").Value = 19.17525512825234

$ws.Range("B20").Value = 12.78965639067267
$ws.Range("C20").Value = 8.619164077622278
$ws.Range("D20").Value = 11.72369081895219
$ws.Range("F20").Value = 27.41633640245032
$ws.Range("G20").Value = 24.33144460712649
$ws.Range("H20").Value = 12.95830461477744
$ws.Range("I20").Value = 17.88067576963624
$ws.Range("J20").Value = 11.14363441023726
$ws.Range("M20").Value = 17.034024404658
$ws.Range("N20").Value = 16.70864467611369
$ws.Range("O20").Value = 19.16515896249361

$ws.Range("B21").Value = 13.39893503766242
$ws.Range("C21").Value = 9.061556735567409
$ws.Range("D21").Value = 11.76258015739517
$ws.Range("F21").Value = 27.44657919350072
$ws.Range("G21").Value = 24.44381658408776
$ws.Range("H21").Value = 12.9231481331422
$ws.Range("I21").Value = 17.78351154004719
$ws.Range("J21").Value = 11.11823087267003
$ws.Range("M21").Value = 17.27199560024549
$ws.Range("N21").Value = 16.67193788672397
$ws.Range("O21").Value = 19.13676511445461

$ws.Range("B22").Value = 13.78242301369991
$ws.Range("C22").Value = 9.338272944521151
$ws.Range("D22").Value = 11.78959583946547
$ws.Range("F22").Value = 27.47263224365168
$ws.Range("G22").Value = 24.52486587078701
$ws.Range("H22").Value = 12.90215437839492
$ws.Range("I22").Value = 17.72336828106204
$ws.Range("J22").Value = 11.10306187135951
$ws.Range("M22").Value = 17.42756330293173
$ws.Range("N22").Value = 16.64937983122674
$ws.Range("O22").Value = 19.12233521592976

$ws.Range("B23").Value = 13.5791143210183
$ws.Range("C23").Value = 9.191724706735354
$ws.Range("D23").Value = 11.77503241379234
$ws.Range("F23").Value = 27.45815344072889
$ws.Range("G23").Value = 24.48091826143631
$ws.Range("H23").Value = 12.9131775741808
$ws.Range("I23").Value = 17.75516011057634
$ws.Range("J23").Value = 11.1110266630578
$ws.Range("M23").Value = 17.34454400579622
$ws.Range("N23").Value = 16.66128863870459
$ws.Range("O23").Value = 19.12965869147229

$ws.Range("B24").Value = 12.77945451089876
$ws.Range("C24").Value = 8.61172500064245
$ws.Range("D24").Value = 11.72308255851366
$ws.Range("F24").Value = 27.41594855447911
$ws.Range("G24").Value = 24.32973764834209
$ws.Range("H24").Value = 12.95891267973836
$ws.Range("I24").Value = 17.88232049217797
$ws.Range("J24").Value = 11.14407381900255
$ws.Range("M24").Value = 17.03013864890226
$ws.Range("N24").Value = 16.70926884275918
$ws.Range("O24").Value = 19.1656926199603

$ws.Range("B25").Value = 11.85611312416225
$ws.Range("C25").Value = 7.933407643657815
$ws.Range("D25").Value = 11.67407314118215
$ws.Range("F25").Value = 27.39753941135528
$ws.Range("G25").Value = 24.19985879898531
$ws.Range("H25").Value = 13.01674649344801
$ws.Range("I25").Value = 18.0338677476272
$ws.Range("J25").Value = 11.17097157113754
$ws.Range("M25").Value = 16.80594192072392
$ws.Range("N25").Value = 16.7468524197598
$ws.Range("O25").Value = 19.20083052657372
